# Update the "想去人数" (interested-people count) figures that changed
# between scrapes, per commit "Update gh-pages to output generated at 456a3b4".

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 5779
$ws1.Range("F3").Value = 13
$ws1.Range("F5").Value = 974
$ws1.Range("F6").Value = 62

# Sheet "演出" (Performance)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 4

# Sheet "全部类型" (All types) - aggregated view of the above sheets
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 5779
$ws4.Range("F3").Value = 13
$ws4.Range("F5").Value = 974
$ws4.Range("F6").Value = 62
$ws4.Range("F7").Value = 4
